$d = $word.ActiveDocument

# Step 1: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("ParentText is a chatbot service that helps you complete your parenting goals using a curriculum designed by Parenting for Lifelong Health with UNICEF and the World Health Organization, and tested all over the world. This programme works! ", $true, $false, $false, $false, $false, $true, 1, $false, "I-ParentText yinkundla yokuxoxa ekusiza ekutheni ufeze izinhloso zakho zokuba umzali ngokusebenzisa ikharikhulamu eyakhiwe yi-Parenting for Lifelong Health ibambisene no-UNICEF kanye ne-World Health Organisation, futhi ehlolwe emhlabeni wonke. Loluhlelo luyasebenza! ", 1)
if (-not $ok) { throw "Step 1 find failed" }

# Step 2: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Being here shows how much you care about providing the best support for your teen. Halala!", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuba khona kwakho lana kutshengisa ukuthi ukukhathalele ukunikeza umntwana wakho ukwesekwa okuvelele. Halala!", 1)
if (-not $ok) { throw "Step 2 find failed" }

# Step 3: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Remember: it is what you do with your teen that will make a difference. ParentText will provide you with tips and skills to help you with your relationship with your teen, but it is up to you to put these tips into practice!", $true, $false, $false, $false, $false, $true, 1, $false, "Khumbula: ilokho okwenza nomntwana wakho okwenza umehluko. I-ParentText izokuhlinzeka ngamacebo namakhono azokusiza ebudlelwaneni bakho nomntwana wakho, kodwa kukuwe ukuwasebenzisa lamacebo!", 1)
if (-not $ok) { throw "Step 3 find failed" }

# Step 4: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("I’m Ayanda, your guide. I may look like a human, but I’m actually a robot produced by Parenting for Lifelong Health and UNICEF to help you learn. ", $true, $false, $false, $false, $false, $true, 1, $false, "Ngingu Ayanda, umhlahlandlela wakho. Ngingabukeka ngathi ngingumuntu, kwodwa ngiyirobhothi elizokusiza eLakhiwe yi-Parenting for Lifelong Health no UNICEF. ", 1)
if (-not $ok) { throw "Step 4 find failed" }

# Step 5: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Today, I’m going to explain how to use ParentText. Together we will review: ", $true, $false, $false, $false, $false, $true, 1, $false, "Namhlanje ngizokuchazela ukuthi isetshenziswa kanjani i-ParentText. Ndawonye sizobuyekeza: ", 1)
if (-not $ok) { throw "Step 5 find failed" }

# Step 6: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("How to earn your Positive Parenting Trophy", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthi uyizuza kanjani iNdondo yokuba uMzali oMuhle", 1)
if (-not $ok) { throw "Step 6 find failed" }

# Step 7: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("How to make progress in your parenting goals", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthi uzithuthukisa kanjani izinhloso zakho zobuzali", 1)
if (-not $ok) { throw "Step 7 find failed" }

# Step 8: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("How to track your progress", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthi uyilandelela kanjani inqubekela phambili yakho", 1)
if (-not $ok) { throw "Step 8 find failed" }

# Step 9: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("How to get help with this course", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthi ulithola kanjani usizo kulesisifundo", 1)
if (-not $ok) { throw "Step 9 find failed" }

# Step 10: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Accessing support to troubleshoot common parenting challenges, and", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthola ukwesekwa ekuxazululeni izinselelo zobuzali ezivamile, ne", 1)
if (-not $ok) { throw "Step 10 find failed" }

# Step 11: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Resources available to you in an emergency or crisis. ", $true, $false, $false, $false, $false, $true, 1, $false, "Izinsiza ongazithola lapho ubhekana nesimo esiphuthumayo noma inhlekele. ", 1)
if (-not $ok) { throw "Step 11 find failed" }

# Step 12: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Earning Your Positive Parenting Trophy", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuthola iNdondo Yakho yokuba uMzali Omuhle", 1)
if (-not $ok) { throw "Step 12 find failed" }

# Step 13: multi-run (break-separated) replace, 3 segments
$rngA = $d.Content
$okA = $rngA.Find.Execute("First, How to complete the course and earn your Positive Parenting Trophy.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okA) { throw "Step 13 find-start failed" }
$segStart = $rngA.Start
$rngB = $d.Content
$okB = $rngB.Find.Execute(" You can choose which goals you want to work on first, but each goal must be completed in order to complete the course and earn your Positive Parenting Trophy.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okB) { throw "Step 13 find-end failed" }
$segEnd = $rngB.End
$full = $d.Range($segStart, $segEnd)
$full.Text = "Okokuqala, Usiqeda kanjani isifundo bese uthola iNdondo Yakho yokuba uMzali Omuhle.`v`v Ngokuncika kwiphrofiyili yakho, lesisifundo sicazwe ngezinghloso zobuzali eziwu 8 noma 9 ezizokusiza ekuthuthukiseni ubudlelwane bakho nomntwana wakho futhi kusize ekutheni umntwana wakho aqhakaze. Inhloso ngayinye ithatha phakathi kwezinsuku ezimbili ukuya kwezine ukuthi uyiqede, usuku nosuku luza nesifundo esisha.`v`v Ungakhetha ukuthi iziphi izinjongo ofuna ukuqala ngazo, kodwa injongo nenjongo kumele iqedwe khona kuzoqedeka isifundo bese uthola iNdondo Yakho yokuba uMzali Omuhle."

# Step 14: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Completing a Goal", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuqeda iNjongo", 1)
if (-not $ok) { throw "Step 14 find failed" }

# Step 15: multi-run (break-separated) replace, 4 segments
$rngA = $d.Content
$okA = $rngA.Find.Execute(" Now, let's learn how to complete each of the goals in the programme.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okA) { throw "Step 15 find-start failed" }
$segStart = $rngA.Start
$rngB = $d.Content
$okB = $rngB.Find.Execute(" You must complete all the skills within a goal to earn a badge. After you have completed one parenting goal, you can select another. Once all of your goal badges are earned, you will have completed the course, and will receive the Positive Parenting Trophy.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okB) { throw "Step 15 find-end failed" }
$segEnd = $rngB.End
$full = $d.Range($segStart, $segEnd)
$full.Text = " Manje masifunde ukuthi siziqeda kanjani lezinjongo ezikuloluhlelo.`v`v Maduze, uzocelwa ukuthi ukhethe inhloso yakho yokuqala. Uma inhloso isikhethiwe, uzoqedela izifundondo khona uzozuza amakhono amasha.`v`v Ikhono elisha litholakala nsukuzonke. Kuthatha ngaphansi kwemizuzu emihlanu ukuqeda iningi lalamakhono. Uma ungakwazanga ukuqeda ikhono owabelwe lona, ngizokubuza ukuthi uyafuna yini ukuliqedela ngakusasa.`v`v Kumele uwaqede wonke amakhono akuleyonhloso yesifundo ukuze uthole ibheji. Uma usuyiqedile inhloso eyodwa yobuzali, ungakhetha enye. Uma usuwathole wonke amabheji enhloso, uzobe ususiqedile isifundo bese uthola iNdondo yakho yoBuzali Obuhle."

# Step 16: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Tracking Progress", $true, $false, $false, $false, $false, $true, 1, $false, "Ukulandelela inqubekelaphambili", 1)
if (-not $ok) { throw "Step 16 find failed" }

# Step 17: multi-run (break-separated) replace, 4 segments
$rngA = $d.Content
$okA = $rngA.Find.Execute("As you move through each day's lesson, you'll receive updates on your progress that look like this: . These check marks tell you how far along you are on the day's lesson.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okA) { throw "Step 17 find-start failed" }
$segStart = $rngA.Start
$rngB = $d.Content
$okB = $rngB.Find.Execute(" Here you can see your progress, review the goals you have achieved, and those that are still incomplete.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okB) { throw "Step 17 find-end failed" }
$segEnd = $rngB.End
$full = $d.Range($segStart, $segEnd)
$full.Text = "Uma uqhubeka nesifundo sosuku, uzothola izibuyekezo ngenqubekelaphambili yakho ebukeka kanje:  Lezizimpawu zikutshela ukuthi usuhambe kanganani esifundweni sakho sosuku.`v`v Uma ufuna ukubona ukuthi usuhambe kangakanani ngezihloso zakho zobuzali, ungabheka kwi- Main Menu. Ukuthola imenyu, bhala u 'Menu`" noma yinini.`v`v Inketho yokuqala kwimenyu imakwe ngokuthi `"bheka inqubekelaphambili yami`".`v`v Kulapha ke lapho ubona inqubekela phambili yakho, ubuyekeze izinhloso ozifezile nalezo ongakazifezi."

# Step 18: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Menu ", $true, $false, $false, $false, $false, $true, 1, $false, "Menyu ", 1)
if (-not $ok) { throw "Step 18 find failed" }

# Step 19: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("The menu contains other features that might help you, too. ", $true, $false, $false, $false, $false, $true, 1, $false, "Imenyu iqukethe ezinye izici ezingakusiza. ", 1)
if (-not $ok) { throw "Step 19 find failed" }

# Step 20: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("In addition to tracking your progress, you can also: ", $true, $false, $false, $false, $false, $true, 1, $false, "Ngokwengeziwe ekubhekeni inqubekelaphambili yakho, unga: ", 1)
if (-not $ok) { throw "Step 20 find failed" }

# Step 21: single-run replace
$rng = $d.Content
$ok = $rng.Find.Execute("Share ParentText with a friend and help them enroll. ", $true, $false, $false, $false, $false, $true, 1, $false, "Yabelana nomngani i-ParentText bese uyamsiza abhalise. ", 1)
if (-not $ok) { throw "Step 21 find failed" }
